$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 296, shifting rows 296:303 down to 297:304
$ws.Rows(296).Insert()

# Populate the new row 296 with the new record (same categorical data as
# neighbouring "Femacal de La Calera" rows, new date/price observations)
$ws.Range("A296").Value = 3
$ws.Range("B296").Value = "Femacal de La Calera"
$ws.Range("C296").Value = "Coquimbo"
$ws.Range("D296").Value = 44939
$ws.Range("E296").Value = 5
$ws.Range("F296").Value = "Fruta"
$ws.Range("G296").Value = 100101
$ws.Range("H296").Value = "Berries"
$ws.Range("I296").Value = 100101001
$ws.Range("J296").Value = "Arándano (blue)"
$ws.Range("K296").Value = "Sin especificar"
$ws.Range("L296").Value = "Primera"
$ws.Range("M296").Value = 62
$ws.Range("N296").Value = 3800
$ws.Range("O296").Value = 4000
$ws.Range("P296").Value = 3903
$ws.Range("Q296").Value = "$/bandeja 2 kilos"
$ws.Range("R296").Value = "Provincia de Curicó"
$ws.Range("S296").Value = 1952
$ws.Range("T296").Value = 2
